$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 9 data
$ws.Range("A9").Value = "Sprint 4"
$ws.Range("B9").Value = "Error message"
$ws.Range("C9").Value = "An error message should be displayed when the server does not answer, or does not answer correctly"
$ws.Range("D9").Value = "Open"

# Copy style from row 8 (A,B,D -> style like col A8/B8/D8, C -> style like C8)
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)  # xlPasteFormats

# Match the taller wrapped-text row height used by similar rows (e.g. row 6)
$ws.Rows.Item(9).RowHeight = 60
$excel.CutCopyMode = $false

# Update view: clear topLeftCell scroll, set selection to E7
$ws.Activate()
$ws.Range("E7").Select()
